# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the zh-cn and de-de handback report sheets (regenerated
# report values).

$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 2 and 3 share the same handoff/handback timestamps
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2:E3").Value = "2016-03-19 12:16:48"
$wsZh.Range("H2:H3").Value = "2016-03-19 12:17:08"

# de-de sheet: rows 2 and 3 share the same handoff/handback timestamps
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2:E3").Value = "2016-03-19 12:16:51"
$wsDe.Range("H2:H3").Value = "2016-03-19 12:17:13"
